# Auto-generated edit script: update cryptos Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.763.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.47%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("E6").Value = "  +0.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3865"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07855"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9898"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.898.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.990"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.18%  "

$ws.Range("E14").Value = "  +0.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06982"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001005"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.71%  "

$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.785.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.281"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.101"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.126.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.831"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.982"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09335"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.09%  "

$ws.Range("E32").Value = "  -1.60%  "

$ws.Range("E33").Value = "  +1.09%  "

$ws.Range("E34").Value = "  +1.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.321"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05785"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.150"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02071"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.663"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5638"
$ws.Range("D40").Style = "Normal"

$ws.Range("E41").Value = "  +1.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.816"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07213"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5298"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.121"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.122"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.834"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.418"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.85%  "

$ws.Range("E51").Value = "  +0.25%  "
